$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("A1").Value2 -eq "Input Sheet") {
        $ws.Range("A1").Value = "Car Name"
    }
    if ($ws.Range("B1").Value2 -eq "Value") {
        $ws.Range("B1").Value = "Values"
    }
}
